# Adding new trigger counter to register spreadsheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bold the header row (row 1), matching the new bold-header styling.
$ws.Range("A1:G1").Font.Bold = $true

# Append the new register row describing the internal trigger counter.
$ws.Range("A39").Value = "internal_trigger_count"
$ws.Range("B39").Value = 32
$ws.Range("C39").Value = "cosmic_ray_internal_trigger_count"
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 32
$ws.Range("F39").Value = "readonly"
$ws.Range("G39").Value = "Counts the number of triggers output by the coincidence_trigger block. These are the internally generated triggers that made it past the veto."

# Keep printable area oriented as portrait (matches saved page setup).
$ws.PageSetup.Orientation = 1

# Restore the view/selection state used when the edit was made.
$ws.Range("E44").Select() | Out-Null
